$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: update Problem Statement text (Problem 1) and wrap it ---
$b3Text = "Sum of elements equal to target in an array.  `nModification : No duplicates, only unique pairs"
$ws.Range("B3").Value = $b3Text
$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 48

# --- Row 4: new Problem 2 entry ---
$ws.Range("A4").Value = 2

$b4Text = "Print element and number of repetitions with the least extra space. Numbers will be repeated in a group and not randomly placed. Every number repeats atleast twice. `nEx : 2 2 3 3 3 3 4 4 5 5 5 -> input 2 2 3 4 4 2 5 3 -> output`nModification : Number may repeat only once`nEx : 2 3 3 4 5 5 -> input  2 1 3 2 4 1 5 2 -> output"
$ws.Range("B4").Value = $b4Text
$ws.Range("B4").WrapText = $true

$ws.Range("C4").Value = "Array manipulation"
$ws.Range("E4").Value = "Java"

$ws.Range("F4").NumberFormat = "mmm-yy"
$ws.Range("F4").Value = 39083

$ws.Rows.Item(4).RowHeight = 110

# --- Column C widened to fit the new "Array manipulation" / "HashMap" labels ---
$ws.Columns.Item(3).ColumnWidth = 15.998697916666666

# --- Final selection left on I4, matching the author's last-saved cursor position ---
$ws.Range("I4").Select() | Out-Null
